$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.388.27'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').Value = '3.303.03'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '525.28'
$ws.Range('E5').Value = '  -2.15%  '
$ws.Range('D6').Value = '172.30'
$ws.Range('E6').Value = '  -6.99%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.344.71'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.583'
$ws.Range('E8').Value = '  -4.32%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '0.602'
$ws.Range('E10').Value = '  -3.99%  '
$ws.Range('D11').Value = '52.57'
$ws.Range('E11').Value = '  -13.97%  '
$ws.Range('D12').Value = '0.133'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '0.0000256'
$ws.Range('E13').Value = '  -4.44%  '
$ws.Range('D14').Value = '8.88'
$ws.Range('E14').Value = '  -3.80%  '
$ws.Range('D15').Value = '3.791.52'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').Value = '3.277.87'
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').Value = '64.229.44'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').Value = '17.35'
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('D20').Value = '11.13'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').Value = '0.953'
$ws.Range('E21').Value = '  -1.97%  '
$ws.Range('D22').Value = '377.86'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  +5.69%  '
$ws.Range('D24').Value = '81.18'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').Value = '11.11'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').Value = '3.68'
$ws.Range('E26').Value = '  -5.10%  '
$ws.Range('D27').Value = '6.17'
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').Value = '2.69'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').Value = '11.18'
$ws.Range('E29').Value = '  -4.88%  '
$ws.Range('D30').Value = '8.07'
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('D31').Value = '28.53'
$ws.Range('E31').Value = '  -2.74%  '
$ws.Range('D32').Value = '626.22'
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('E33').Value = '  -5.70%  '
$ws.Range('D34').Value = '11.14'
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('D35').Value = '0.105'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').Value = '56.54'
$ws.Range('E36').Value = '  -5.61%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = '35.78'
$ws.Range('E38').Value = '  -3.56%  '
$ws.Range('D39').Value = '0.376'
$ws.Range('E39').Value = '  -5.55%  '
$ws.Range('D40').Value = '0.0₃0735'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').Value = '0.989'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').Value = '3.16'
$ws.Range('E42').Value = '  +8.13%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.125'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').Value = '2.860.42'
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('D46').Value = '2.69'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').Value = '0.0395'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  -5.40%  '
$ws.Range('D49').Value = '3.02'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = '137.53'
$ws.Range('E50').Value = '  +1.48%  '
$ws.Range('D51').Value = '0.124'
$ws.Range('E51').Value = '  -2.90%  '
